$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns D and E
$ws.Range("D1").Value = "Update Date"
$ws.Range("E1").Value = "Comment"

# Column widths for new columns D and E (target stored widths: 12.85546875 / 37.42578125;
# the host's ColumnWidth setter quantizes to 1/6-character steps, so we feed the input
# that lands closest to the target after that quantization)
$ws.Columns.Item(4).ColumnWidth = 12.0
$ws.Columns.Item(5).ColumnWidth = 36.666666666666664

# Row 54: DHC_RX - remove C54, add D54 (date) and E54 (comment)
$ws.Range("C54").ClearContents()
$ws.Range("D54").Value = 43544
$ws.Range("D54").NumberFormat = "mm-dd-yy"
$ws.Range("E54").Value = "Unassign from Pod Digestive Health"

# Row 56: ENDO - add C56, D56 (date), E56 (comment)
$ws.Range("C56").Value = "Digestive Health"
$ws.Range("D56").Value = 43544
$ws.Range("E56").Value = "Assign to Pod Digestive Health"

# Reuse the same date style (xf) for D56 as D54 instead of minting a new one
$ws.Range("D54").Copy()
$ws.Range("D56").PasteSpecial(-4122)
